$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new "Ngày thực hiện" / "{{Exports.Date}}" column (N) ---------------
# Write the values first, then clone formatting from the neighbouring columns
# that already carry the styles we want (L1/M1 header style, M2 value style).
# This re-uses the existing style/font/fill/border records instead of minting
# new duplicate ones.
$ws.Range("N1").Value = "Ngày thực hiện"
$ws.Range("N2").Value = "{{Exports.Date}}"

$ws.Range("L1").Copy() | Out-Null
$ws.Range("N1").PasteSpecial(-4122) | Out-Null

$ws.Range("M2").Copy() | Out-Null
$ws.Range("N2").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Resize / split the columns around the new reward-date column -------------
# Columns I (9) and J (10) used to share the same width as E:J; give them their
# own (narrower) widths so the reward code/name/value/date block reads better.
$ws.Columns("I").ColumnWidth = 13.418
$ws.Columns("J").ColumnWidth = 12.586
# Mã số giải thưởng / Giải thưởng columns get narrower too.
$ws.Columns("K").ColumnWidth = 21.836
$ws.Columns("L").ColumnWidth = 20.17
# New "Ngày thực hiện" column width.
$ws.Columns("N").ColumnWidth = 16.67

# --- Row 2 (data row) grows taller to fit the extra wrapped column ------------
$ws.Rows("2").RowHeight = 45

# --- Selection, as left by the author after the edit ---------------------------
$ws.Range("C5").Select() | Out-Null
